{"js": "// Replace the date line and every two-digit x two-digit multiplication\n// problem in the practice sheet with the values from the new session.\nconst replacements = [\n  [\"2024-12-03 Tuesday\", \"2024-12-04 Wednesday\"],\n  [\"71\u00d759=\", \"40\u00d796=\"],\n  [\"92\u00d757=\", \"86\u00d715=\"],\n  [\"98\u00d783=\", \"77\u00d711=\"],\n  [\"44\u00d770=\", \"38\u00d725=\"],\n  [\"74\u00d772=\", \"85\u00d719=\"],\n  [\"73\u00d789=\", \"44\u00d766=\"],\n  [\"70\u00d716=\", \"55\u00d789=\"],\n  [\"41\u00d730=\", \"99\u00d735=\"],\n  [\"27\u00d729=\", \"35\u00d712=\"],\n  [\"86\u00d726=\", \"61\u00d734=\"],\n  [\"82\u00d772=\", \"99\u00d773=\"],\n  [\"64\u00d765=\", \"43\u00d751=\"],\n  [\"85\u00d767=\", \"76\u00d729=\"],\n  [\"88\u00d789=\", \"91\u00d739=\"],\n  [\"49\u00d747=\", \"25\u00d776=\"],\n  [\"94\u00d761=\", \"92\u00d743=\"],\n  [\"86\u00d740=\", \"52\u00d718=\"],\n  [\"84\u00d741=\", \"70\u00d713=\"],\n  [\"74\u00d786=\", \"58\u00d752=\"],\n  [\"60\u00d772=\", \"73\u00d793=\"],\n  [\"40\u00d739=\", \"16\u00d735=\"],\n  [\"52\u00d717=\", \"73\u00d754=\"],\n  [\"77\u00d790=\", \"17\u00d725=\"],\n  [\"33\u00d726=\", \"60\u00d748=\"],\n  [\"96\u00d791=\", \"27\u00d757=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2024-12-03 Tuesday\", \"2024-12-04 Wednesday\"),\n  @(\"71\u00d759=\", \"40\u00d796=\"),\n  @(\"92\u00d757=\", \"86\u00d715=\"),\n  @(\"98\u00d783=\", \"77\u00d711=\"),\n  @(\"44\u00d770=\", \"38\u00d725=\"),\n  @(\"74\u00d772=\", \"85\u00d719=\"),\n  @(\"73\u00d789=\", \"44\u00d766=\"),\n  @(\"70\u00d716=\", \"55\u00d789=\"),\n  @(\"41\u00d730=\", \"99\u00d735=\"),\n  @(\"27\u00d729=\", \"35\u00d712=\"),\n  @(\"86\u00d726=\", \"61\u00d734=\"),\n  @(\"82\u00d772=\", \"99\u00d773=\"),\n  @(\"64\u00d765=\", \"43\u00d751=\"),\n  @(\"85\u00d767=\", \"76\u00d729=\"),\n  @(\"88\u00d789=\", \"91\u00d739=\"),\n  @(\"49\u00d747=\", \"25\u00d776=\"),\n  @(\"94\u00d761=\", \"92\u00d743=\"),\n  @(\"86\u00d740=\", \"52\u00d718=\"),\n  @(\"84\u00d741=\", \"70\u00d713=\"),\n  @(\"74\u00d786=\", \"58\u00d752=\"),\n  @(\"60\u00d772=\", \"73\u00d793=\"),\n  @(\"40\u00d739=\", \"16\u00d735=\"),\n  @(\"52\u00d717=\", \"73\u00d754=\"),\n  @(\"77\u00d790=\", \"17\u00d725=\"),\n  @(\"33\u00d726=\", \"60\u00d748=\"),\n  @(\"96\u00d791=\", \"27\u00d757=\"),\n)\n\nforeach ($pair in $replacements) {\n  $rng = $d.Content\n  $find = $rng.Find\n  $find.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $pair[1]\n  $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2) | Out-Null\n}\n"}
